# "Improved Get and Post, still working on Patch"
# Adds a per-endpoint checklist (GetById / GetAll / GetByRangeOfId / Patch /
# Delete) for the Get/Post-ish rows, moves the status-code notes out to
# column P, notes "Heavily changed controller" next to Automapper, and
# switches the active tab back to Sheet1.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# --- Sheet1: clear the cells that moved / disappeared -----------------
$ws1.Range("E1").ClearContents()
$ws1.Range("E3").ClearContents()
$ws1.Range("J4").ClearContents()
$ws1.Range("J5").ClearContents()
$ws1.Range("J6").ClearContents()
$ws1.Range("J7").ClearContents()

# --- Row 9 (Automapper): controller note (written early, matches the
#     author's original edit order) ------------------------------------
$ws1.Range("L9").Value = "Heavily changed controller"

# --- Row 2: new header row for the per-endpoint columns ---------------
$ws1.Range("C2").Value = "GetById"
$ws1.Range("D2").Value = "GetAll"
$ws1.Range("E2").Value = "GetByRangeOfId"
$ws1.Range("F2").Value = "Patch"
$ws1.Range("G2").Value = "Delete"
$ws1.Range("K2").Value = "Optional"

# --- Row 4 (Repository): mark all endpoints done, move note to K ------
$ws1.Range("C4").Value = "v"
$ws1.Range("D4").Value = "v"
$ws1.Range("E4").Value = "v"
$ws1.Range("F4").Value = "v"
$ws1.Range("G4").Value = "v"
$ws1.Range("K4").Value = "When delete fail, show related entity"

# --- Row 5 (RepoLogging): status-code note moves to column P ----------
$ws1.Range("P5").Value = "NotFound()"

# --- Row 6 (Ofm Get): mark endpoints done, status note moves to P -----
$ws1.Range("C6").Value = "v"
$ws1.Range("D6").Value = "v"
$ws1.Range("E6").Value = "v"
$ws1.Range("F6").Value = "v"
$ws1.Range("G6").Value = "v"
$ws1.Range("P6").Value = "BadRequest(modelState)"

# --- Row 7 (Ofm Post): mark endpoints done, status note moves to P ----
$ws1.Range("C7").Value = "v"
$ws1.Range("D7").Value = "v"
$ws1.Range("E7").Value = "v"
$ws1.Range("F7").Value = "v"
$ws1.Range("G7").Value = "v"
$ws1.Range("P7").Value = "NoContent()"

# --- Row 8 (Ofm Patch): mark endpoints done, status note moves to P ---
$ws1.Range("C8").Value = "v"
$ws1.Range("D8").Value = "v"
$ws1.Range("E8").Value = "v"
$ws1.Range("F8").Value = "v"
$ws1.Range("G8").Value = "v"
$ws1.Range("P8").Value = "Ok()"

# --- Column widths for the new endpoint columns (best fit) ------------
$ws1.Columns.Item(3).ColumnWidth = 7.1666666666666667
$ws1.Columns.Item(4).ColumnWidth = 5.7369791666666667
$ws1.Columns.Item(5).ColumnWidth = 14.8776041666666666
$ws1.Columns.Item(6).ColumnWidth = 5.0221354166666667
$ws1.Columns.Item(7).ColumnWidth = 6.1666666666666667

# --- Switch the active tab back to Sheet1 ------------------------------
$ws1.Activate() | Out-Null
$ws1.Range("H6").Select() | Out-Null
